$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (ALC)
$ws.Range("H6").Value = 375.66666
$ws.Range("I6").Value = 97.07692
$ws.Range("J6").Value = 1100
$ws.Range("K6").Value = 291.23076
$ws.Range("L6").Value = 3300
$ws.Range("M6").Value = -179.23076
$ws.Range("N6").Value = -3524

# Row 8 (ALC)
$ws.Range("H8").Value = 11.833333
$ws.Range("I8").Value = 11.833333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 35.499999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 103.500001

# Row 107 (ALC)
$ws.Range("H107").Value = 634.4091
$ws.Range("I107").Value = 696.3
$ws.Range("J107").Value = 582.8333
$ws.Range("K107").Value = 696.3
$ws.Range("L107").Value = 582.8333
$ws.Range("M107").Value = 1223.7
$ws.Range("N107").Value = -4422.8333

# Row 112 (ALC)
$ws.Range("H112").Value = 1728
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1728
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 5184
$ws.Range("N112").Value = -7400

# Row 113 (ALC)
$ws.Range("H113").Value = 4445.3516
$ws.Range("I113").Value = 3834.261
$ws.Range("J113").Value = 5449.2856
$ws.Range("K113").Value = 3834.261
$ws.Range("L113").Value = 5449.2856
$ws.Range("M113").Value = -580.261
$ws.Range("N113").Value = -11957.2856

# Row 129 (ALC)
$ws.Range("H129").Value = 781.3125
$ws.Range("I129").Value = 432.16666
$ws.Range("J129").Value = 990.8
$ws.Range("K129").Value = 1296.49998
$ws.Range("L129").Value = 2972.4
$ws.Range("M129").Value = 3703.50002
$ws.Range("N129").Value = -12972.4

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1317.1904
$ws.Range("I2").Value = 1401.091
$ws.Range("J2").Value = 1224.9
$ws.Range("K2").Value = 1401.091
$ws.Range("L2").Value = 1224.9
$ws.Range("M2").Value = -1288.091
$ws.Range("N2").Value = -1450.9

# Row 110 (ARM)
$ws.Range("H110").Value = 1379.2858
$ws.Range("I110").Value = 1237.2727
$ws.Range("J110").Value = 1900
$ws.Range("K110").Value = 1237.2727
$ws.Range("L110").Value = 1900
$ws.Range("M110").Value = 807.7273

# Row 116 (ARM)
$ws.Range("H116").Value = 1317.1904
$ws.Range("I116").Value = 1401.091
$ws.Range("J116").Value = 1224.9
$ws.Range("K116").Value = 1401.091
$ws.Range("L116").Value = 1224.9
$ws.Range("M116").Value = 892.9090000000001
$ws.Range("N116").Value = -5812.9

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1317.1904
$ws.Range("I3").Value = 1401.091
$ws.Range("J3").Value = 1224.9
$ws.Range("K3").Value = 1401.091
$ws.Range("L3").Value = 1224.9
$ws.Range("M3").Value = -1287.091
$ws.Range("N3").Value = -1452.9

# Row 87 (BSM)
$ws.Range("H87").Value = 35400
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 35400
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 35400
$ws.Range("N87").Value = -37896

# Row 90 (BSM)
$ws.Range("H90").Value = 35400
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 35400
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 106200
$ws.Range("N90").Value = -118680

# Row 94 (BSM)
$ws.Range("H94").Value = 1134.0769
$ws.Range("I94").Value = 843.7143
$ws.Range("J94").Value = 3674.75
$ws.Range("K94").Value = 843.7143
$ws.Range("L94").Value = 3674.75
$ws.Range("M94").Value = -392.7143
$ws.Range("N94").Value = -4576.75

# Row 107 (BSM)
$ws.Range("H107").Value = 2101.8333
$ws.Range("I107").Value = 2101.8333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2101.8333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -181.8332999999998
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (CRP)
$ws.Range("H4").Value = 3059.8235
$ws.Range("I4").Value = 2553.3333
$ws.Range("J4").Value = 3168.3572
$ws.Range("K4").Value = 2553.3333
$ws.Range("L4").Value = 3168.3572
$ws.Range("M4").Value = -2441.3333
$ws.Range("N4").Value = -3392.3572

# Row 16 (CRP)
$ws.Range("H16").Value = 1007.5
$ws.Range("I16").Value = 650
$ws.Range("J16").Value = 1126.6666
$ws.Range("K16").Value = 650
$ws.Range("L16").Value = 1126.6666
$ws.Range("M16").Value = -363
$ws.Range("N16").Value = -1700.6666

# Row 39 (CRP)
$ws.Range("H39").Value = 6461.25
$ws.Range("I39").Value = 6461.25
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6461.25
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -6070.25

# Row 49 (CRP)
$ws.Range("H49").Value = 6461.25
$ws.Range("I49").Value = 6461.25
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 6461.25
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -6279.25

# Row 107 (CRP)
$ws.Range("H107").Value = 983.4706
$ws.Range("I107").Value = 613.2353000000001
$ws.Range("J107").Value = 1353.7059
$ws.Range("K107").Value = 613.2353000000001
$ws.Range("L107").Value = 1353.7059
$ws.Range("M107").Value = 1306.7647
$ws.Range("N107").Value = -5193.7059

# Row 113 (CRP)
$ws.Range("H113").Value = 1007.5
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 1126.6666
$ws.Range("K113").Value = 650
$ws.Range("L113").Value = 1126.6666
$ws.Range("M113").Value = 1520
$ws.Range("N113").Value = -5466.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (CUL)
$ws.Range("H4").Value = 1192.1538
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 1318
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 3954
$ws.Range("M4").Value = -1388
$ws.Range("N4").Value = -4178

# Row 20 (CUL)
$ws.Range("H20").Value = 3180
$ws.Range("I20").Value = 800
$ws.Range("J20").Value = 4766.6665
$ws.Range("K20").Value = 2400
$ws.Range("L20").Value = 14299.9995
$ws.Range("M20").Value = -2173
$ws.Range("N20").Value = -14753.9995

# Row 131 (CUL)
$ws.Range("H131").Value = 1080.3928
$ws.Range("I131").Value = 775.25
$ws.Range("J131").Value = 1095.65
$ws.Range("K131").Value = 2325.75
$ws.Range("L131").Value = 3286.95
$ws.Range("M131").Value = 2714.25
$ws.Range("N131").Value = -13366.95

# Row 133 (CUL)
$ws.Range("H133").Value = 7700
$ws.Range("I133").Value = 2000
$ws.Range("J133").Value = 8333.333000000001
$ws.Range("K133").Value = 6000
$ws.Range("L133").Value = 24999.999
$ws.Range("M133").Value = -940
$ws.Range("N133").Value = -35119.999

# Row 138 (CUL)
$ws.Range("H138").Value = 1228.8235
$ws.Range("I138").Value = 1092.6666
$ws.Range("J138").Value = 2250
$ws.Range("K138").Value = 3277.9998
$ws.Range("L138").Value = 6750
$ws.Range("M138").Value = 1862.0002
$ws.Range("N138").Value = -17030

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (GSM)
$ws.Range("H113").Value = 1661.8182
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1811.4286
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 1811.4286
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -6151.4286

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (LTW)
$ws.Range("H2").Value = 216543
$ws.Range("I2").Value = 362950.25
$ws.Range("J2").Value = 21333.334
$ws.Range("K2").Value = 362950.25
$ws.Range("L2").Value = 21333.334
$ws.Range("M2").Value = -362838.25
$ws.Range("N2").Value = -21557.334

# Row 13 (LTW)
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# Row 61 (LTW)
$ws.Range("H61").Value = 1401.5
$ws.Range("I61").Value = 946.5454999999999
$ws.Range("J61").Value = 2402.4
$ws.Range("K61").Value = 946.5454999999999
$ws.Range("L61").Value = 2402.4
$ws.Range("M61").Value = -744.5454999999999
$ws.Range("N61").Value = -2806.4

# Row 113 (LTW)
$ws.Range("H113").Value = 1401.5
$ws.Range("I113").Value = 946.5454999999999
$ws.Range("J113").Value = 2402.4
$ws.Range("K113").Value = 946.5454999999999
$ws.Range("L113").Value = 2402.4
$ws.Range("M113").Value = 1223.4545
$ws.Range("N113").Value = -6742.4

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (WVR)
$ws.Range("H2").Value = 17833.334
$ws.Range("I2").Value = 90000
$ws.Range("J2").Value = 3400
$ws.Range("K2").Value = 90000
$ws.Range("L2").Value = 3400
$ws.Range("M2").Value = -89888
$ws.Range("N2").Value = -3624

# Row 86 (WVR)
$ws.Range("H86").Value = 24728.25
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 24728.25
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 24728.25
$ws.Range("N86").Value = -26974.25

# Row 87 (WVR)
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 89 (WVR)
$ws.Range("H89").Value = 24728.25
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 24728.25
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 123641.25
$ws.Range("N89").Value = -134873.25

# Row 90 (WVR)
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 109 (WVR)
$ws.Range("H109").Value = 30000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 30000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
